# Added test classes for triangle builders -> reflect new "double" boundary
# test rows on the "Invalid" data sheet (min/max double values), and make
# that sheet the active one, like the authoring Excel session ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invalid")

# Grow the existing table (Таблица2) to cover the two new data rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C10"))

# Labels first (matches shared-string insertion order of the source edit).
$ws.Range("D9").Value = "double_Max_value"
$ws.Range("D10").Value = "double_Min_value"

# Row 9: double.MaxValue as text in each of the three side columns.
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "1,7976931348623157E+308"
$ws.Range("A9").Font.Color = 2763306

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "1,7976931348623157E+308"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "1,7976931348623157E+308"

# Row 10: double.MinValue (negative max) as text in each side column.
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "-1,7976931348623157E+308"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "-1,7976931348623157E+308"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "-1,7976931348623157E+308"

# Match the selection left in the saved file and make "Invalid" the active tab.
$ws.Range("E10").Select()
$ws.Activate()
